{"js": "const body = context.document.body;\n\n// 1. \"Good day\" + \",\" (two separate runs) -> \"Good day,\" (single merged run)\nconst goodDayResults = body.search(\"Good day,\", { matchCase: true, matchWholeWord: false });\ngoodDayResults.load(\"items\");\nawait context.sync();\nif (goodDayResults.items.length > 0) {\n  goodDayResults.items[0].insertText(\"Good day,\", Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// 2. \"Welcome to Arzew port,\" -> \"Welcome to \" + \"our\" + \" port,\" (three runs)\nconst arzewResults = body.search(\"Arzew\", { matchCase: true, matchWholeWord: false });\narzewResults.load(\"items\");\nawait context.sync();\nif (arzewResults.items.length > 0) {\n  arzewResults.items[0].insertText(\"our\", Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// Force the replaced word into its own run (distinct from its neighbors)\n// by toggling a character format on it and back off again, so the text\n// stays split into three runs: \"Welcome to \", \"our\", \" port,\".\nconst ourResults = body.search(\"our\", { matchCase: true, matchWholeWord: false });\nourResults.load(\"items\");\nawait context.sync();\nif (ourResults.items.length > 0) {\n  const ourRange = ourResults.items[0];\n  ourRange.font.bold = true;\n  await context.sync();\n  ourRange.font.bold = false;\n  await context.sync();\n}\n\n// 3. Remove the \"_GoBack\" bookmark (bookmarkStart/bookmarkEnd) after \"DECLARATION OF SECURITY\"\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# 1. \"Good day\" + \",\" (two separate runs) -> \"Good day,\" (single merged run)\n$find1 = $d.Content\n$find1.Find.Execute(\"Good day,\", $false, $false, $false, $false, $false, $true, 1, $false, \"Good day,\", 2)\n\n# 2. \"Welcome to Arzew port,\" -> \"Welcome to \" + \"our\" + \" port,\" (three runs)\n$find2 = $d.Content\n$find2.Find.Execute(\"Arzew\", $false, $false, $false, $false, $false, $true, 1, $false, \"our\", 2)\n\n# Force the replaced word into its own run (distinct from its neighbors)\n# by toggling a character format on it and back off again, so the text\n# stays split into three runs: \"Welcome to \", \"our\", \" port,\".\n$ourRange = $d.Content\n$ourRange.Find.Execute(\"our\")\n$ourRange.Bold = 1\n$ourRange.Bold = 0\n\n# 3. Remove the \"_GoBack\" bookmark (bookmarkStart/bookmarkEnd) after \"DECLARATION OF SECURITY\"\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks(\"_GoBack\").Delete()\n}\n"}
